$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Title) values, row by row, so new shared strings are
# appended in the same order the target workbook uses. ---
$ws.Range("B2").Value = "Anadolu notları"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "1984"
$ws.Range("B3").ClearFormats()

$ws.Range("A4").Value = 3
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B4").Value = "1703 Edirne vakası"

$ws.Range("A5").Value = 4
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("B5").Value = "Halide Edib Adivar"

$ws.Range("A6").Value = 5
$ws.Range("A2").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("B6").Value = "Karmasik Duygular"

$ws.Range("A7").Value = 6
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("B7").Value = "Tarihte ilginç gerçekler"

# --- Column C (Authors) values, row by row. ---
$ws.Range("C2").Value = "['Reşat Nuri Güntekin']"
$ws.Range("C3").Value = "['George Orwell']"
$ws.Range("C4").Value = "['Tahir Sevinç', 'Behset Karaca', 'Süleyman Demirel Üniversitesi. Sosyal Bilimler Enstitüsü. Tarih Anabilim Dalı']"
$ws.Range("C5").Value = "['Sinekli Bakkal']"
$ws.Range("C6").Value = "['Stefan Zweig']"
$ws.Range("C7").Value = "['']"
